$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names, URLs) ---
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

# --- Numeric-looking text cells (price/volume columns) ---
# Force these to remain plain text (matching original inlineStr cells)
# by temporarily applying a text number format, then clearing the format
# afterwards so no extra styling is left on the cell.
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '61.042.98'
$r.ClearFormats()
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '  -1.86%  '
$r.ClearFormats()
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '2.973.02'
$r.ClearFormats()
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '  -0.45%  '
$r.ClearFormats()
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '  +0.09%  '
$r.ClearFormats()
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '595.36'
$r.ClearFormats()
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '  +2.60%  '
$r.ClearFormats()
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '142.01'
$r.ClearFormats()
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '  -2.45%  '
$r.ClearFormats()
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '  +0.08%  '
$r.ClearFormats()
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.514'
$r.ClearFormats()
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '  -1.26%  '
$r.ClearFormats()
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '2.972.32'
$r.ClearFormats()
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '  -0.43%  '
$r.ClearFormats()
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '6.00'
$r.ClearFormats()
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '  +6.42%  '
$r.ClearFormats()
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.145'
$r.ClearFormats()
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '  -1.94%  '
$r.ClearFormats()
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '  +2.96%  '
$r.ClearFormats()
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '  -0.45%  '
$r.ClearFormats()
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '33.95'
$r.ClearFormats()
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '  -1.57%  '
$r.ClearFormats()
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '  +2.11%  '
$r.ClearFormats()
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '3.464.67'
$r.ClearFormats()
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '  -0.16%  '
$r.ClearFormats()
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '61.042.29'
$r.ClearFormats()
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '  -1.78%  '
$r.ClearFormats()
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '6.85'
$r.ClearFormats()
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '  -2.30%  '
$r.ClearFormats()
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '2.969.65'
$r.ClearFormats()
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = '  -0.49%  '
$r.ClearFormats()
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '446.48'
$r.ClearFormats()
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '  -1.97%  '
$r.ClearFormats()
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '14.06'
$r.ClearFormats()
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '  +1.75%  '
$r.ClearFormats()
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '0.678'
$r.ClearFormats()
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '  +0.29%  '
$r.ClearFormats()
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '  -0.14%  '
$r.ClearFormats()
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '81.87'
$r.ClearFormats()
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '  +2.61%  '
$r.ClearFormats()
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '  -5.11%  '
$r.ClearFormats()
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '11.88'
$r.ClearFormats()
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '  -2.46%  '
$r.ClearFormats()
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '10.26'
$r.ClearFormats()
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = '  +2.98%  '
$r.ClearFormats()
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = '  +0.08%  '
$r.ClearFormats()
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = '  +2.74%  '
$r.ClearFormats()
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = '  +0.07%  '
$r.ClearFormats()
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '7.06'
$r.ClearFormats()
$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = '  -2.22%  '
$r.ClearFormats()
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = '  -2.13%  '
$r.ClearFormats()
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '27.02'
$r.ClearFormats()
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = '  +1.16%  '
$r.ClearFormats()
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '0.107'
$r.ClearFormats()
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = '  +0.99%  '
$r.ClearFormats()
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = '  +3.59%  '
$r.ClearFormats()
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = '  -1.27%  '
$r.ClearFormats()
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '5.73'
$r.ClearFormats()
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = '  +0.32%  '
$r.ClearFormats()
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '50.16'
$r.ClearFormats()
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = '  +0.49%  '
$r.ClearFormats()
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '  -2.93%  '
$r.ClearFormats()
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '8.97'
$r.ClearFormats()
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '  +0.23%  '
$r.ClearFormats()
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '  +9.12%  '
$r.ClearFormats()
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '  -2.64%  '
$r.ClearFormats()
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '389.69'
$r.ClearFormats()
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '  -4.62%  '
$r.ClearFormats()
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '38.94'
$r.ClearFormats()
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '  +1.97%  '
$r.ClearFormats()
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '  -0.68%  '
$r.ClearFormats()
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '  -4.02%  '
$r.ClearFormats()
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '2.674.73'
$r.ClearFormats()
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '  -3.04%  '
$r.ClearFormats()
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '130.24'
$r.ClearFormats()
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '  +2.43%  '
$r.ClearFormats()
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '  +0.07%  '
$r.ClearFormats()
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.106'
$r.ClearFormats()
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '  -0.83%  '
$r.ClearFormats()
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '  -0.62%  '
$r.ClearFormats()
